$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $oCell = $ws.Cells.Item($r, 15)
    $pCell = $ws.Cells.Item($r, 16)
    $oVal = $oCell.Value2
    $pVal = $pCell.Value2
    $oCell.Value2 = $pVal
    $pCell.Value2 = $oVal
}
